# Add two new worksheets, Figure_40 and Figure_41, at the end of the workbook,
# each containing the same county/race-ethnicity table layout used by the
# other Figure_* sheets (header row + 7 data rows, columns A-H).

$wb = $excel.ActiveWorkbook

$headers = @("county", "Asian/Pacific Islander", "Black", "Hispanic/Latino", "Multiracial/Other", "Native American", "White", "Total")

$figure40Rows = @(
    @("Imperial", 3.54, 14.23, 16.08, 9.800000000000001, 8.73, 5.06, 14.56),
    @("Los Angeles", 5.48, 10.1, 6.8, 9.1, 6.77, 6.77, 6.9),
    @("Orange", 5.39, 8.59, 5.44, 8.289999999999999, 5.73, 5.13, 5.45),
    @("Riverside", 5.36, 10.63, 7.74, 8.02, 10.14, 6.44, 7.35),
    @("San Bernardino", 5.89, 11.66, 7.74, 10.44, 11.79, 6.67, 7.68),
    @("Ventura", 4.25, 9.59, 6.13, 6.39, 6.83, 5.08, 5.58),
    @("SCAG", 5.44, 10.29, 6.96, 8.83, 8.029999999999999, 6.22, 6.78)
)

$figure41Rows = @(
    @("Imperial", 24.75, 38.35, 32.03, 2.11, 45.86, 13.1, 29.69),
    @("Los Angeles", 13.79, 18.54, 28.32, 14.32, 16.23, 9.93, 19.9),
    @("Orange", 11.76, 16.08, 24.26, 9.41, 15.61, 7.23, 14.09),
    @("Riverside", 13.08, 16.19, 25.05, 15.15, 21.4, 10.6, 18.56),
    @("San Bernardino", 15.76, 23.46, 24.91, 18, 27.77, 12.64, 20.49),
    @("Ventura", 6.76, 12.5, 24.21, 9.44, 2.83, 7.34, 14.58),
    @("SCAG", 13.17, 18.7, 26.84, 13.58, 18.41, 9.470000000000001, 18.6)
)

function Add-FigureSheet($sheetName, $dataRows) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $ws.Name = $sheetName

    for ($col = 1; $col -le $headers.Length; $col++) {
        $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
    }

    $headerRange = $ws.Range("A1:H1")
    $headerRange.Font.Bold = $true
    $headerRange.HorizontalAlignment = -4108

    $rowNum = 2
    foreach ($row in $dataRows) {
        for ($col = 1; $col -le $row.Length; $col++) {
            $ws.Cells.Item($rowNum, $col).Value = $row[$col - 1]
        }
        $rowNum++
    }
}

Add-FigureSheet "Figure_40" $figure40Rows
Add-FigureSheet "Figure_41" $figure41Rows
